$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.056.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.651.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.64%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5272"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2598"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.74%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06311"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.29"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07786"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.515"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.84%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.671.26"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.97%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.878.47"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5472"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8180"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.34"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.060.31"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.61%  "
$ws.Range("E19").Value = "  -0.03%  "
$ws.Range("E20").Value = "  -0.67%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "190.26"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.20%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.005"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("E24").Value = "  -0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1231"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.208"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.93%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.98"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.23%  "
$ws.Range("E29").Value = "  +0.84%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05807"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("E32").Value = "  -0.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.256"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.592"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.80%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.795"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("E36").Value = "  -0.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9423"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5748"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01606"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8489"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.10%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.55"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.98%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -4.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.028.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.793.06"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "56.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.68%  "
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("E48").Value = "  -0.34%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.847"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05137"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.446"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.49%  "
